$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.278.03"
$ws.Range("E2").Value = "  +1.53%  "

$ws.Range("D3").Value = "1.806.58"
$ws.Range("E3").Value = "  +3.30%  "

$ws.Range("E4").Value = "  -0.15%  "

$ws.Range("D5").Value = "'338.36"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.64%  "

$ws.Range("E6").Value = "  -0.08%  "

$ws.Range("D7").Value = "'0.4645"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +19.89%  "

$ws.Range("D8").Value = "'0.3826"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +13.21%  "

$ws.Range("D9").Value = "'45.42"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.83%  "

$ws.Range("D10").Value = "'1.161"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.26%  "

$ws.Range("E11").Value = "  +5.64%  "

$ws.Range("E12").Value = "  +0.22%  "

$ws.Range("E13").Value = "  -0.11%  "

$ws.Range("D14").Value = "'6.371"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.13%  "

$ws.Range("D15").Value = "'7.516"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +6.20%  "

$ws.Range("D16").Value = "1.813.30"
$ws.Range("E16").Value = "  +3.63%  "

$ws.Range("D17").Value = "'0.00001096"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.48%  "

$ws.Range("D18").Value = "'0.06738"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.05%  "

$ws.Range("D19").Value = "'81.68"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.04%  "

$ws.Range("D20").Value = "'1.000"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.13%  "

$ws.Range("D21").Value = "'17.65"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +5.24%  "

$ws.Range("D22").Value = "'6.440"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.19%  "

$ws.Range("D23").Value = "28.276.01"
$ws.Range("E23").Value = "  +1.45%  "

$ws.Range("D24").Value = "'11.92"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.40%  "

$ws.Range("E25").Value = "  +1.20%  "

$ws.Range("D26").Value = "'20.73"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.53%  "

$ws.Range("D27").Value = "'154.03"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.04%  "

$ws.Range("D28").Value = "'2.375"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.29%  "

$ws.Range("D29").Value = "2.013.35"
$ws.Range("E29").Value = "  +3.25%  "

$ws.Range("D30").Value = "'133.00"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.55%  "

$ws.Range("D31").Value = "'1.253"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.29%  "

$ws.Range("E32").Value = "  +0.39%  "

$ws.Range("D33").Value = "'0.09631"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +9.50%  "

$ws.Range("D34").Value = "'5.864"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.64%  "

$ws.Range("D35").Value = "'0.2352"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +12.09%  "

$ws.Range("B36").Value = "InternetComputer(DFINITY)"
$ws.Range("C36").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D36").Value = "'5.292"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.23%  "

$ws.Range("D37").Value = "'12.12"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.07%  "

$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value = "'0.02359"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.60%  "

$ws.Range("B39").Value = "Hedera"
$ws.Range("C39").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D39").Value = "'0.06367"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.00%  "

$ws.Range("D40").Value = "'0.6636"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.62%  "

$ws.Range("D41").Value = "'1.247"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.47%  "

$ws.Range("B42").Value = "WEMIXTOKEN"
$ws.Range("C42").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D42").Value = "'1.492"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.79%  "

$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").Value = "'8.366"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.46%  "

$ws.Range("E44").Value = "  +4.12%  "

$ws.Range("D45").Value = "'1.001"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.06%  "

$ws.Range("D46").Value = "'0.6151"
$ws.Range("D46").Style = "Normal"

$ws.Range("D47").Value = "'3.859"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.27%  "

$ws.Range("D48").Value = "'131.76"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.86%  "

$ws.Range("D49").Value = "'2.056"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.18%  "

$ws.Range("D50").Value = "'1.182"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.77%  "

$ws.Range("D51").Value = "'0.07163"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.84%  "
